$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Set the "Title" value (row 5, column B) to match the Name value.
$ws.Range("B5").Value = "DroitExerciceComplementaire"

# Update the "Date" value (row 8, column B) to the new timestamp.
$ws.Range("B8").Value = "2025-07-17T14:35:50+00:00"
